# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new F value for the "展览" sheet
$exhibitionUpdates = @{
    5  = 5089
    7  = 47
    9  = 563
    10 = 520
    13 = 1418
    14 = 3757
    17 = 131
    18 = 87
    19 = 2744
    20 = 138
    21 = 34
    25 = 72
    26 = 16
    27 = 131
    28 = 64
    29 = 280
    30 = 48
    31 = 4
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new F value for the "全部类型" sheet
$allTypesUpdates = @{
    6  = 5089
    8  = 47
    10 = 563
    11 = 520
    14 = 1418
    15 = 3757
    18 = 131
    19 = 87
    20 = 2744
    21 = 138
    22 = 34
    26 = 72
    27 = 16
    28 = 131
    29 = 64
    30 = 280
    31 = 48
    32 = 4
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
